$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. The "_GoBack" bookmark currently sits at the end of the paragraph
#    that ends with "...to ensure that it is correct." Remove it from
#    there - it needs to move to the start of the "Stability testing"
#    bullet further down.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 2. Insert a brand-new justified paragraph right after the "Also, I
#    have obviously not verified..." paragraph, introducing the sample
#    output file. Build it via InsertXML so the run layout (including
#    the spell-check proofErr markers around "outputfile" and
#    "birdypigs") matches exactly what Word itself would produce.
# ---------------------------------------------------------------------
$anchorPara = $d.Paragraphs.Item(9)
$anchorPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(10)
$newRange = $newPara.Range
$newRange.End = $newRange.End - 1   # exclude the paragraph mark

$newParaXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">A small sample output is included in the file called </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>outputfile</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>birdypigs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> directory.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$newRange.InsertXML($newParaXml)

# ---------------------------------------------------------------------
# 3. Re-add the "_GoBack" bookmark at the very beginning of the
#    "Stability testing" bullet, which - now that the new paragraph has
#    been inserted above it - is paragraph 11.
# ---------------------------------------------------------------------
$stabilityPara = $d.Paragraphs.Item(11)
$bookmarkSpot = $d.Range($stabilityPara.Range.Start, $stabilityPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)
